$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: "LoadingDetail" / "System (A)" rows ---
# Header cell G7 should look like the rest of row 7 (copy format from F7)
$ws.Range("F7").Copy() | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null
$ws.Range("G7").Value = "LoadingDetail"

# Data cells G8/G9 should look like the rest of their rows (copy format from E8/E9)
$ws.Range("E8").Copy() | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null
$ws.Range("G8").Value = "System (A)"

$ws.Range("E9").Copy() | Out-Null
$ws.Range("G9").PasteSpecial(-4122) | Out-Null
$ws.Range("G9").Value = "System (A)"

$excel.CutCopyMode = 0

# --- Update the selection shown when the sheet is reopened ---
$ws.Range("G7:G9").Select() | Out-Null
